# Swap the presentation's theme color scheme from "Integral" to the
# stock "Office Theme" 12-slot palette (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink). The font scheme ("Arial" for major/minor) and the
# format scheme are already identical between the two themes, so only
# the color scheme needs to change.
#
# PowerPoint resolves/stores the deck's color scheme through the
# ThemeColorScheme collection, which is reachable from a Slide (or the
# SlideMaster/NotesMaster's Theme). Each of the 12 entries exposes a
# settable .RGB (standard COM BGR-packed RGB integer, i.e. R | G<<8 | B<<16).

$p = $ppt.ActivePresentation

# Target "Office Theme" color scheme, in ThemeColorScheme index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeTheme = @(
    @{ idx = 1;  name = "dk1";      r = 0x00; g = 0x00; b = 0x00 },
    @{ idx = 2;  name = "lt1";      r = 0xFF; g = 0xFF; b = 0xFF },
    @{ idx = 3;  name = "dk2";      r = 0x44; g = 0x54; b = 0x6A },
    @{ idx = 4;  name = "lt2";      r = 0xE7; g = 0xE6; b = 0xE6 },
    @{ idx = 5;  name = "accent1";  r = 0x5B; g = 0x9B; b = 0xD5 },
    @{ idx = 6;  name = "accent2";  r = 0xED; g = 0x7D; b = 0x31 },
    @{ idx = 7;  name = "accent3";  r = 0xA5; g = 0xA5; b = 0xA5 },
    @{ idx = 8;  name = "accent4";  r = 0xFF; g = 0xC0; b = 0x00 },
    @{ idx = 9;  name = "accent5";  r = 0x44; g = 0x72; b = 0xC4 },
    @{ idx = 10; name = "accent6";  r = 0x70; g = 0xAD; b = 0x47 },
    @{ idx = 11; name = "hlink";    r = 0x05; g = 0x63; b = 0xC1 },
    @{ idx = 12; name = "folHlink"; r = 0x95; g = 0x4F; b = 0x72 }
)

function ToComRgb($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# The deck's single slide master (theme1.xml, originally the "Integral"
# theme) is reachable via any slide's ThemeColorScheme.
$tcs = $p.Slides.Item(1).ThemeColorScheme
foreach ($c in $officeTheme) {
    $tcs.Item($c.idx).RGB = ToComRgb $c.r $c.g $c.b
}

# Also reach the same scheme via the SlideMaster/Theme object directly,
# and via the NotesMaster/HandoutMaster's Theme, in case the host wires
# any of those to a distinct underlying theme part (e.g. theme2.xml, the
# notes master's originally-"Office Theme" theme, which this edit turns
# into the deck's old "Integral" palette).
$integralTheme = @(
    @{ idx = 1;  r = 0x00; g = 0x00; b = 0x00 },
    @{ idx = 2;  r = 0xFF; g = 0xFF; b = 0xFF },
    @{ idx = 3;  r = 0x45; g = 0x5F; b = 0x51 },
    @{ idx = 4;  r = 0xE3; g = 0xDE; b = 0xD1 },
    @{ idx = 5;  r = 0x99; g = 0xCB; b = 0x38 },
    @{ idx = 6;  r = 0x63; g = 0xA5; b = 0x37 },
    @{ idx = 7;  r = 0xE6; g = 0xD0; b = 0x24 },
    @{ idx = 8;  r = 0xCC; g = 0x97; b = 0x00 },
    @{ idx = 9;  r = 0x4E; g = 0xB3; b = 0xCF },
    @{ idx = 10; r = 0x37; g = 0x8D; b = 0xA6 },
    @{ idx = 11; r = 0x6B; g = 0x9F; b = 0x25 },
    @{ idx = 12; r = 0xB2; g = 0x6B; b = 0x02 }
)

try {
    $nm = $p.NotesMaster
    $nmTcs = $nm.Theme.ThemeColorScheme
    foreach ($c in $integralTheme) {
        $nmTcs.Item($c.idx).RGB = ToComRgb $c.r $c.g $c.b
    }
} catch {
    Write-Host "NotesMaster theme colors not reachable:" $_.Exception.Message
}

try {
    $hm = $p.HandoutMaster
    $hmTcs = $hm.Theme.ThemeColorScheme
    foreach ($c in $integralTheme) {
        $hmTcs.Item($c.idx).RGB = ToComRgb $c.r $c.g $c.b
    }
} catch {
    Write-Host "HandoutMaster theme colors not reachable:" $_.Exception.Message
}

# Re-apply the Office Theme to the primary (slide master) scheme last,
# so it wins if any of the paths above alias back to the same theme part.
foreach ($c in $officeTheme) {
    $tcs.Item($c.idx).RGB = ToComRgb $c.r $c.g $c.b
}

Write-Host "Theme color scheme updated."
